$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new "zoom user ids" sheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "zoom user ids"

$zoomUserIds = @(
    "iclrconf+Walcott-Bryant@gmail.com",
    "iclrconf+Kaelbling@gmail.com",
    "iclrconf+Benjamin@gmail.com",
    "iclrconf+Dinh@gmail.com",
    "iclrconf+Schaar@gmail.com",
    "iclrconf+Parikh@gmail.com",
    "iclrconf+LeCunBengio@gmail.com",
    "iclrconf+Jordan@gmail.com"
)

for ($i = 0; $i -lt $zoomUserIds.Length; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $zoomUserIds[$i]
}

$ws2.Range("A1:A8").Select() | Out-Null

# --- Update Sheet1: rename the zoom_user_id column to host_zoom_user_email
#     and point every row at the shared Zoom-host account ---
$ws1.Range("K1").Value = "host_zoom_user_email"
$ws1.Range("K2:K9").Value = "leetncamp+iclr@gmail.com"

# Martha White's panelist contact changed to her new ICLR-issued Zoom address
$ws1.Range("J6").Value = "Martha White <iclrconf+mwhite@gmail.com>"

$ws1.Activate() | Out-Null
$ws1.Range("K2:K9").Select() | Out-Null
